$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (20 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 12162.5
$ws.Range("J43").Value = 11716.667
$ws.Range("L43").Value = 11716.667
$ws.Range("N43").Value = -11854.667
$ws.Range("H70").Value = 12487.308
$ws.Range("J70").Value = 6875
$ws.Range("L70").Value = 20625
$ws.Range("N70").Value = -21165
$ws.Range("H73").Value = 12487.308
$ws.Range("J73").Value = 6875
$ws.Range("L73").Value = 20625
$ws.Range("N73").Value = -22497
$ws.Range("H132").Value = 3495.0232
$ws.Range("I132").Value = 1579.75
$ws.Range("K132").Value = 4739.25
$ws.Range("M132").Value = -2209.25
$ws.Range("H138").Value = 2037.3208
$ws.Range("J138").Value = 2592.9656
$ws.Range("L138").Value = 7778.8968
$ws.Range("N138").Value = -18058.8968

# --- Sheet: ARM (34 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 13265.5
$ws.Range("I28").Value = 7514.3335
$ws.Range("K28").Value = 7514.3335
$ws.Range("M28").Value = -7322.3335
$ws.Range("H41").Value = 2625.5
$ws.Range("I41").Value = 2625.5
$ws.Range("K41").Value = 2625.5
$ws.Range("M41").Value = -2211.5
$ws.Range("H74").Value = 2112.25
$ws.Range("I74").Value = 1885.7391
$ws.Range("J74").Value = 3154.2
$ws.Range("K74").Value = 1885.7391
$ws.Range("L74").Value = 3154.2
$ws.Range("M74").Value = -1011.7391
$ws.Range("N74").Value = -4902.2
$ws.Range("H77").Value = 2112.25
$ws.Range("I77").Value = 1885.7391
$ws.Range("J77").Value = 3154.2
$ws.Range("K77").Value = 9428.6955
$ws.Range("L77").Value = 15771
$ws.Range("M77").Value = -5060.6955
$ws.Range("N77").Value = -24507
$ws.Range("H80").Value = 19995
$ws.Range("J80").Value = 19995
$ws.Range("L80").Value = 19995
$ws.Range("N80").Value = -21991
$ws.Range("H83").Value = 19995
$ws.Range("J83").Value = 19995
$ws.Range("L83").Value = 59985
$ws.Range("N83").Value = -69969
$ws.Range("H99").Value = 13265.5
$ws.Range("I99").Value = 7514.3335
$ws.Range("K99").Value = 7514.3335
$ws.Range("M99").Value = -4519.3335

# --- Sheet: BSM (12 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1278.5333
$ws.Range("I80").Value = 930.1429000000001
$ws.Range("K80").Value = 930.1429000000001
$ws.Range("M80").Value = 67.85709999999995
$ws.Range("H83").Value = 1278.5333
$ws.Range("I83").Value = 930.1429000000001
$ws.Range("K83").Value = 4650.7145
$ws.Range("M83").Value = 341.2855
$ws.Range("H105").Value = 2706207.8
$ws.Range("I105").Value = 3128614
$ws.Range("K105").Value = 3128614
$ws.Range("M105").Value = -3126867

# --- Sheet: CRP (25 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 902.3333
$ws.Range("I19").Value = 183
$ws.Range("J19").Value = 4499
$ws.Range("K19").Value = 183
$ws.Range("L19").Value = 4499
$ws.Range("M19").Value = -13
$ws.Range("N19").Value = -4839
$ws.Range("H24").Value = 902.3333
$ws.Range("I24").Value = 183
$ws.Range("J24").Value = 4499
$ws.Range("K24").Value = 183
$ws.Range("L24").Value = 4499
$ws.Range("M24").Value = -13
$ws.Range("N24").Value = -4839
$ws.Range("H107").Value = 1738.9048
$ws.Range("I107").Value = 1532.2
$ws.Range("J107").Value = 1926.8182
$ws.Range("K107").Value = 1532.2
$ws.Range("L107").Value = 1926.8182
$ws.Range("M107").Value = 387.8
$ws.Range("N107").Value = -5766.8182
$ws.Range("H132").Value = 3527.9
$ws.Range("I132").Value = 3422.875
$ws.Range("K132").Value = 10268.625
$ws.Range("M132").Value = -7738.625

# --- Sheet: CUL (39 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 545.2727
$ws.Range("J11").Value = 999
$ws.Range("L11").Value = 2997
$ws.Range("N11").Value = -3277
$ws.Range("H39").Value = 2619.8
$ws.Range("I39").Value = 774.75
$ws.Range("J39").Value = 10000
$ws.Range("K39").Value = 2324.25
$ws.Range("L39").Value = 30000
$ws.Range("M39").Value = -2030.25
$ws.Range("N39").Value = -30588
$ws.Range("H46").Value = 317.2
$ws.Range("I46").Value = 339
$ws.Range("J46").Value = 230
$ws.Range("K46").Value = 1017
$ws.Range("L46").Value = 690
$ws.Range("M46").Value = -926
$ws.Range("N46").Value = -872
$ws.Range("H62").Value = 9749.25
$ws.Range("I62").Value = 7999
$ws.Range("J62").Value = 11499.5
$ws.Range("K62").Value = 23997
$ws.Range("L62").Value = 34498.5
$ws.Range("M62").Value = -23311
$ws.Range("N62").Value = -35870.5
$ws.Range("H65").Value = 9749.25
$ws.Range("I65").Value = 7999
$ws.Range("J65").Value = 11499.5
$ws.Range("K65").Value = 71991
$ws.Range("L65").Value = 103495.5
$ws.Range("M65").Value = -68559
$ws.Range("N65").Value = -110359.5
$ws.Range("H116").Value = 3406.1667
$ws.Range("I116").Value = 2087.8
$ws.Range("J116").Value = 9998
$ws.Range("K116").Value = 6263.400000000001
$ws.Range("L116").Value = 29994
$ws.Range("M116").Value = -2821.400000000001
$ws.Range("N116").Value = -36878

# --- Sheet: GSM (43 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 34876252
$ws.Range("I11").Value = 19315430
$ws.Range("K11").Value = 19315430
$ws.Range("M11").Value = -19315291
$ws.Range("H21").Value = 29740.666
$ws.Range("I21").Value = 22000
$ws.Range("J21").Value = 33611
$ws.Range("K21").Value = 22000
$ws.Range("L21").Value = 33611
$ws.Range("M21").Value = -21827
$ws.Range("N21").Value = -33957
$ws.Range("H24").Value = 40010670
$ws.Range("I24").Value = 200000000
$ws.Range("J24").Value = 13337.5
$ws.Range("K24").Value = 200000000
$ws.Range("L24").Value = 13337.5
$ws.Range("M24").Value = -199999827
$ws.Range("N24").Value = -13683.5
$ws.Range("H30").Value = 29740.666
$ws.Range("I30").Value = 22000
$ws.Range("J30").Value = 33611
$ws.Range("K30").Value = 22000
$ws.Range("L30").Value = 33611
$ws.Range("M30").Value = -21895
$ws.Range("N30").Value = -33821
$ws.Range("H99").Value = 18710
$ws.Range("I99").Value = 14437.571
$ws.Range("J99").Value = 28679
$ws.Range("K99").Value = 14437.571
$ws.Range("L99").Value = 28679
$ws.Range("M99").Value = -12191.571
$ws.Range("N99").Value = -33171
$ws.Range("H102").Value = 2472.889
$ws.Range("I102").Value = 2034.7142
$ws.Range("J102").Value = 4006.5
$ws.Range("K102").Value = 2034.7142
$ws.Range("L102").Value = 4006.5
$ws.Range("M102").Value = -412.7141999999999
$ws.Range("N102").Value = -7250.5
$ws.Range("H132").Value = 3783
$ws.Range("I132").Value = 3783
$ws.Range("K132").Value = 11349
$ws.Range("M132").Value = -8819

# --- Sheet: LTW (24 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 18491.334
$ws.Range("I23").Value = 18491.334
$ws.Range("K23").Value = 18491.334
$ws.Range("M23").Value = -18261.334
$ws.Range("H33").Value = 9403.75
$ws.Range("I33").Value = 9403.75
$ws.Range("K33").Value = 9403.75
$ws.Range("M33").Value = -9113.75
$ws.Range("H62").Value = 23249
$ws.Range("J62").Value = 23249
$ws.Range("L62").Value = 23249
$ws.Range("N62").Value = -24497
$ws.Range("H65").Value = 23249
$ws.Range("J65").Value = 23249
$ws.Range("L65").Value = 69747
$ws.Range("N65").Value = -75987
$ws.Range("H100").Value = 3157.2778
$ws.Range("I100").Value = 2285
$ws.Range("K100").Value = 2285
$ws.Range("M100").Value = -1744
$ws.Range("H136").Value = 2614.1
$ws.Range("I136").Value = 1975.9231
$ws.Range("K136").Value = 5927.7693
$ws.Range("M136").Value = -3377.7693

# --- Sheet: WVR (11 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 45987.25
$ws.Range("J15").Value = 45987.25
$ws.Range("L15").Value = 45987.25
$ws.Range("N15").Value = -46563.25
$ws.Range("H136").Value = 1003.7826
$ws.Range("I136").Value = 954.3158
$ws.Range("J136").Value = 1238.75
$ws.Range("K136").Value = 2862.9474
$ws.Range("L136").Value = 3716.25
$ws.Range("M136").Value = -312.9474
$ws.Range("N136").Value = -8816.25

Write-Host "Applied 208 cell updates across 8 sheets"